$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# Row 33 (ALC)
$ws1.Range("H33").Value = 42991.332
$ws1.Range("J33").Value = 7333.3335
$ws1.Range("L33").Value = 7333.3335
$ws1.Range("N33").Value = -7791.3335

# Row 38 (ALC)
$ws1.Range("H38").Value = 1698449.4
$ws1.Range("I38").Value = 2304245.8
$ws1.Range("J38").Value = 2220
$ws1.Range("K38").Value = 6912737.399999999
$ws1.Range("L38").Value = 6660
$ws1.Range("M38").Value = -6912365.399999999
$ws1.Range("N38").Value = -7404

# Row 40 (ALC)
$ws1.Range("H40").Value = 40464.77
$ws1.Range("I40").Value = 68559.266
$ws1.Range("J40").Value = 2154.0908
$ws1.Range("K40").Value = 68559.266
$ws1.Range("L40").Value = 2154.0908
$ws1.Range("M40").Value = -68384.266
$ws1.Range("N40").Value = -2504.0908

# Row 62 (ALC)
$ws1.Range("H62").Value = 924.2857
$ws1.Range("I62").Value = 924.2857
$ws1.Range("J62").Value = 0
$ws1.Range("K62").Value = 924.2857
$ws1.Range("L62").Value = 0
$ws1.Range("M62").Value = -300.2857
$ws1.Range("N62").ClearContents()

# Row 65 (ALC)
$ws1.Range("H65").Value = 924.2857
$ws1.Range("I65").Value = 924.2857
$ws1.Range("J65").Value = 0
$ws1.Range("K65").Value = 4621.4285
$ws1.Range("L65").Value = 0
$ws1.Range("M65").Value = -1501.4285
$ws1.Range("N65").Value = 0

# Row 80 (ALC)
$ws1.Range("H80").Value = 38873.848
$ws1.Range("I80").Value = 213.33333
$ws1.Range("J80").Value = 91592.73
$ws1.Range("K80").Value = 639.99999
$ws1.Range("L80").Value = 274778.19
$ws1.Range("M80").Value = 358.00001
$ws1.Range("N80").Value = -276774.19

# Row 83 (ALC)
$ws1.Range("H83").Value = 38873.848
$ws1.Range("I83").Value = 213.33333
$ws1.Range("J83").Value = 91592.73
$ws1.Range("K83").Value = 1919.99997
$ws1.Range("L83").Value = 824334.5699999999
$ws1.Range("M83").Value = 3072.00003
$ws1.Range("N83").Value = -834318.5699999999

# Row 107 (ALC)
$ws1.Range("H107").Value = 502.72223
$ws1.Range("I107").Value = 471.4375
$ws1.Range("K107").Value = 471.4375
$ws1.Range("M107").Value = 1448.5625

# Row 138 (ALC)
$ws1.Range("H138").Value = 6444.231
$ws1.Range("I138").Value = 1282.2059
$ws1.Range("J138").Value = 16194.723
$ws1.Range("K138").Value = 3846.6177
$ws1.Range("L138").Value = 48584.169
$ws1.Range("M138").Value = 1293.3823
$ws1.Range("N138").Value = -58864.169

# Row 45 (ARM)
$ws2.Range("H45").Value = 38261.703
$ws2.Range("I45").Value = 53578.316
$ws2.Range("J45").Value = 1884.75
$ws2.Range("K45").Value = 53578.316
$ws2.Range("L45").Value = 1884.75
$ws2.Range("M45").Value = -53201.316
$ws2.Range("N45").Value = -2638.75

# Row 110 (ARM)
$ws2.Range("H110").Value = 100201144
$ws2.Range("I110").Value = 143143600
$ws2.Range("J110").Value = 2060
$ws2.Range("K110").Value = 143143600
$ws2.Range("L110").Value = 2060
$ws2.Range("M110").Value = -143141555
$ws2.Range("N110").Value = -6150

# Row 107 (BSM)
$ws3.Range("H107").Value = 66724290
$ws3.Range("I107").Value = 125107410
$ws3.Range("J107").Value = 725.2857
$ws3.Range("K107").Value = 125107410
$ws3.Range("L107").Value = 725.2857
$ws3.Range("M107").Value = -125105490
$ws3.Range("N107").Value = -4565.2857

# Row 134 (BSM)
$ws3.Range("H134").Value = 2696.532
$ws3.Range("I134").Value = 2739.5715
$ws3.Range("J134").Value = 2571
$ws3.Range("K134").Value = 8218.7145
$ws3.Range("L134").Value = 7713
$ws3.Range("M134").Value = -5683.7145
$ws3.Range("N134").Value = -12783

# Row 16 (CRP)
$ws4.Range("H16").Value = 1505.5555
$ws4.Range("I16").Value = 1550
$ws4.Range("J16").Value = 1470
$ws4.Range("K16").Value = 1550
$ws4.Range("L16").Value = 1470
$ws4.Range("M16").Value = -1263
$ws4.Range("N16").Value = -2044

# Row 22 (CRP)
$ws4.Range("H22").Value = 470.4
$ws4.Range("I22").Value = 369.5
$ws4.Range("J22").Value = 874
$ws4.Range("K22").Value = 369.5
$ws4.Range("L22").Value = 874
$ws4.Range("M22").Value = -19.5
$ws4.Range("N22").Value = -1574

# Row 96 (CRP)
$ws4.Range("H96").Value = 28000
$ws4.Range("J96").Value = 28000
$ws4.Range("L96").Value = 28000
$ws4.Range("N96").Value = -33492

# Row 107 (CRP)
$ws4.Range("H107").Value = 8350
$ws4.Range("I107").Value = 13349.875
$ws4.Range("J107").Value = 350.2
$ws4.Range("K107").Value = 13349.875
$ws4.Range("L107").Value = 350.2
$ws4.Range("M107").Value = -11429.875
$ws4.Range("N107").Value = -4190.2

# Row 113 (CRP)
$ws4.Range("H113").Value = 1505.5555
$ws4.Range("I113").Value = 1550
$ws4.Range("J113").Value = 1470
$ws4.Range("K113").Value = 1550
$ws4.Range("L113").Value = 1470
$ws4.Range("M113").Value = 620
$ws4.Range("N113").Value = -5810

# Row 132 (CRP)
$ws4.Range("H132").Value = 26788570
$ws4.Range("I132").Value = 25002636
$ws4.Range("J132").Value = 31253406
$ws4.Range("K132").Value = 75007908
$ws4.Range("L132").Value = 93760218
$ws4.Range("M132").Value = -75005378
$ws4.Range("N132").Value = -93765278

# Row 82 (CUL)
$ws5.Range("H82").Value = 2717.5
$ws5.Range("I82").Value = 2246.6667
$ws5.Range("K82").Value = 6740.000100000001
$ws5.Range("M82").Value = -6334.000100000001

# Row 85 (CUL)
$ws5.Range("H85").Value = 2717.5
$ws5.Range("I85").Value = 2246.6667
$ws5.Range("K85").Value = 6740.000100000001
$ws5.Range("M85").Value = -5336.000100000001

# Row 114 (CUL)
$ws5.Range("H114").Value = 739.4545000000001
$ws5.Range("I114").Value = 325.375
$ws5.Range("J114").Value = 1843.6666
$ws5.Range("K114").Value = 976.125
$ws5.Range("L114").Value = 5530.9998
$ws5.Range("M114").Value = 2277.875
$ws5.Range("N114").Value = -12038.9998

# Row 63 (GSM)
$ws6.Range("H63").Value = 24900
$ws6.Range("J63").Value = 24900
$ws6.Range("L63").Value = 24900
$ws6.Range("N63").Value = -26272

# Row 66 (GSM)
$ws6.Range("H66").Value = 24900
$ws6.Range("J66").Value = 24900
$ws6.Range("L66").Value = 74700
$ws6.Range("N66").Value = -81564

# Row 80 (GSM)
$ws6.Range("H80").Value = 91004460
$ws6.Range("I80").Value = 166839820
$ws6.Range("K80").Value = 166839820
$ws6.Range("M80").Value = -166838822

# Row 83 (GSM)
$ws6.Range("H83").Value = 91004460
$ws6.Range("I83").Value = 166839820
$ws6.Range("K83").Value = 834199100
$ws6.Range("M83").Value = -834194108

# Row 107 (GSM)
$ws6.Range("H107").Value = 632112.5600000001
$ws6.Range("I107").Value = 562.5
$ws6.Range("J107").Value = 1263662.6
$ws6.Range("K107").Value = 562.5
$ws6.Range("L107").Value = 1263662.6
$ws6.Range("M107").Value = 1357.5
$ws6.Range("N107").Value = -1267502.6

# Row 113 (GSM)
$ws6.Range("H113").Value = 1894.4445
$ws6.Range("I113").Value = 1625
$ws6.Range("J113").Value = 1971.4286
$ws6.Range("K113").Value = 1625
$ws6.Range("L113").Value = 1971.4286
$ws6.Range("M113").Value = 545
$ws6.Range("N113").Value = -6311.4286

# Row 22 (LTW)
$ws7.Range("H22").Value = 1037.0454
$ws7.Range("J22").Value = 1069.75
$ws7.Range("L22").Value = 1069.75
$ws7.Range("N22").Value = -1659.75

# Row 27 (LTW)
$ws7.Range("H27").Value = 1037.0454
$ws7.Range("J27").Value = 1069.75
$ws7.Range("L27").Value = 1069.75
$ws7.Range("N27").Value = -1283.75

# Row 136 (LTW)
$ws7.Range("H136").Value = 1783.3043
$ws7.Range("I136").Value = 1549.5
$ws7.Range("J136").Value = 2625
$ws7.Range("K136").Value = 4648.5
$ws7.Range("L136").Value = 7875
$ws7.Range("M136").Value = -2098.5
$ws7.Range("N136").Value = -12975

# Row 21 (WVR)
$ws8.Range("H21").Value = 41008.5
$ws8.Range("J21").Value = 41008.5
$ws8.Range("L21").Value = 41008.5
$ws8.Range("N21").Value = -41478.5

# Row 24 (WVR)
$ws8.Range("H24").Value = 353200
$ws8.Range("J24").Value = 353200
$ws8.Range("L24").Value = 353200
$ws8.Range("N24").Value = -353660

# Row 28 (WVR)
$ws8.Range("H28").Value = 15040
$ws8.Range("I28").Value = 10000
$ws8.Range("J28").Value = 16300
$ws8.Range("K28").Value = 10000
$ws8.Range("L28").Value = 16300
$ws8.Range("M28").Value = -9652
$ws8.Range("N28").Value = -16996

# Row 30 (WVR)
$ws8.Range("H30").Value = 35005
$ws8.Range("J30").Value = 10001
$ws8.Range("L30").Value = 10001
$ws8.Range("N30").Value = -10215

# Row 35 (WVR)
$ws8.Range("H35").Value = 41008.5
$ws8.Range("J35").Value = 41008.5
$ws8.Range("L35").Value = 41008.5
$ws8.Range("N35").Value = -41588.5

# Row 64 (WVR)
$ws8.Range("H64").Value = 0
$ws8.Range("J64").Value = 0
$ws8.Range("L64").Value = 0
$ws8.Range("N64").ClearContents()

# Row 67 (WVR)
$ws8.Range("H67").Value = 0
$ws8.Range("J67").Value = 0
$ws8.Range("L67").Value = 0
$ws8.Range("N67").ClearContents()

# Row 81 (WVR)
$ws8.Range("H81").Value = 223317.67
$ws8.Range("J81").Value = 201331.8
$ws8.Range("L81").Value = 402663.6
$ws8.Range("N81").Value = -404785.6

# Row 84 (WVR)
$ws8.Range("H84").Value = 223317.67
$ws8.Range("J84").Value = 201331.8
$ws8.Range("L84").Value = 2013318
$ws8.Range("N84").Value = -2023926

# Row 107 (WVR)
$ws8.Range("H107").Value = 67465.60000000001
$ws8.Range("I107").Value = 600
$ws8.Range("J107").Value = 72241.71000000001
$ws8.Range("K107").Value = 1800
$ws8.Range("L107").Value = 216725.13
$ws8.Range("M107").Value = 120
$ws8.Range("N107").Value = -220565.13

# Row 132 (WVR)
$ws8.Range("H132").Value = 2226.6226
$ws8.Range("I132").Value = 2046.8096
$ws8.Range("J132").Value = 2913.182
$ws8.Range("K132").Value = 6140.4288
$ws8.Range("L132").Value = 8739.545999999998
$ws8.Range("M132").Value = -3610.4288
$ws8.Range("N132").Value = -13799.546

# Row 136 (WVR)
$ws8.Range("H136").Value = 840.8
$ws8.Range("I136").Value = 523.0741
$ws8.Range("K136").Value = 1569.2223
$ws8.Range("M136").Value = 980.7776999999999
